# Add a new localization string "strChkWindowPosition" to the translation
# table, keeping the table's existing alphabetical order (sorted by the
# "Key" column).
#
# The new key sorts between "strChkExportIntegration" (row 7) and
# "strDifferentiationAlgorithms" (old row 8), so a row is physically
# inserted at row 8 - this correctly shifts every following row (and its
# row-height / formatting) down by one, which is simpler and more robust
# than re-running the table's Sort.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a blank worksheet row right above the old row 8, then grow the
# table definition so it covers the extra row at the bottom (B2:P35 ->
# B2:P36).
$ws.Rows.Item(8).Insert()
$tbl.Resize($ws.Range("B2:P36"))

# Fill in the new entry.
$ws.Range("B8").Value = "strChkWindowPosition"
$ws.Range("C8").Value = 'In "settings" form, tab "User interface"'
$ws.Range("D8").Value = "Remember window position and size on startup"

# Match formatting used elsewhere in the table:
#  - Key column (B) uses left/center alignment with word-wrap.
#  - Comment/English columns (C/D) use left/center alignment, word-wrap,
#    and an unlocked protection state (matching other long/wrapped entries
#    such as the one now at D9).
$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").WrapText = $true

$ws.Range("D9").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
